# Faster gate drive on pfc
# Swap the row-11 resistor part for a Viking Tech precision SMD resistor,
# tighten the order quantities, update unit/order price, and mark the
# changed cells with wrap-text formatting. Also update the active-cell
# selection left in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new part data (LCSC C217732 / Viking Tech ARG03FTC0130) ---
$ws.Range("A11").Value = "C217732"
$ws.Range("B11").Value = "ARG03FTC0130"
$ws.Range("C11").Value = "Viking Tech"
$ws.Range("F11").Value = "13 ±1% 1/10W 0603 High Precision & Low TCR SMD Resistors RoHS"
$ws.Range("L11").Value = "https://lcsc.com/product-detail/High-Precision-Low-TCR-SMD-Resistors_Viking-Tech-ARG03FTC0130_C217732.html"

# Order / min-mult quantities drop from 100 to 50
$ws.Range("H11").Value = 50
$ws.Range("I11").Value = 50

# Unit price / order price updated for new part
$ws.Range("J11").Value = 0.0119
$ws.Range("K11").Value = 0.6

# Wrap text on the cells that now carry the new (longer) content
$ws.Range("B11").WrapText = $true
$ws.Range("C11").WrapText = $true
$ws.Range("J11").WrapText = $true
$ws.Range("K11").WrapText = $true

# Row grew slightly taller to fit the wrapped text
$ws.Range("A11:L11").RowHeight = 14.9

# --- Sheet view: move the remembered selection ---
$ws.Range("A11").Select()
